$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Sheet "Statistikk": G2 58 -> 59
# ---------------------------------------------------------------------------
$wsStat = $wb.Worksheets.Item("Statistikk")
$wsStat.Range("G2").Value = 59

# ---------------------------------------------------------------------------
# 2) Sheet "Statistikk": row 14 (digital) D14 2->3, E14 4->5
# ---------------------------------------------------------------------------
$wsStat.Range("D14").Value = 3
$wsStat.Range("E14").Value = 5

# ---------------------------------------------------------------------------
# 3) Sheet "digital": insert a new row 5, pushing the EMPE2500/Bygningssimulering
#    entry (currently in row 4, columns G:I) down into row 5, and replace the
#    row-4 G:I entry with the new EMFE1000 / Matematikk 1000 entry.
# ---------------------------------------------------------------------------
$wsDigital = $wb.Worksheets.Item("digital")

# Grab the values presently sitting in G4:I4 before we move anything.
$oldCode = $wsDigital.Range("G4").Value2
$oldName = $wsDigital.Range("H4").Value2
$oldDesc = $wsDigital.Range("I4").Value2

# Push row 4 downward starting a new row 5 (matches the cell style below it).
$wsDigital.Rows.Item(5).Insert()

# Move the previous EMPE2500 / Bygningssimulering entry into the new row 5.
$wsDigital.Range("G5").Value = $oldCode
$wsDigital.Range("H5").Value = $oldName
$wsDigital.Range("I5").Value = $oldDesc

# Write the new Matematikk 1000 entry into row 4.
$wsDigital.Range("G4").Value = "EMFE1000"
$wsDigital.Range("H4").Value = "Matematikk 1000"
$wsDigital.Range("I4").Value = "Studenten kan vurdere resultater fra matematiske beregninger forklare og bruke grunnleggende numeriske algoritmer som inneholder kodeelementene tilordning for- og while-løkker og if-tester skrive presise forklaringer og begrunnelser til framgangsmåter og demonstrere korrekt bruk av matematisk notasjon vurdere egne og andre studenters faglige arbeider og formulere skriftlige og muntlige vurderinger av disse arbeidene på en faglig korrekt og presis måte verføre et praktisk problem fra eget fagområde til matematisk form slik at det kan løses - analytisk eller numerisk bruke matematiske metoder og digitale verktøy som er relevante for eget fagfelt bruke matematikk til å kommunisere om ingeniørfaglige problemstillinger"

# ---------------------------------------------------------------------------
# 4) Sheet "modell": update the F6 learning-outcome description text.
# ---------------------------------------------------------------------------
$wsModell = $wb.Worksheets.Item("modell")
$wsModell.Range("F6").Value = "Studenten kan anvende den deriverte til å modellere og analysere dynamiske systemer diskutere hvordan ideen bak definisjonen av det bestemte integralet kan brukes til å sette opp integraler for beregning av størrelser drøfte ideene bak noen analytiske og numeriske metoder som brukes for å løse differensiallikninger og sette opp og løse differensiallikninger for praktiske problemer som er relevante innen eget fagområde drøfte metoder for å løse lineære likningssystemer ved hjelp av matriseregning og drøfte numeriske metoder for å løse likninger og sette opp og løse likninger for praktiske problemer fra eget fagområde"
